# Update LDLC prices history
# Insert a new timestamp column ("2026-01-28 15:21:05") right before the
# existing "nom" / "url_produit" columns (which shift one column to the
# right, V -> W and W -> X). For rows that already carry a price in the
# last existing timestamp column (U), the new column repeats that same
# price (mirrors how this scraper's history sheet is appended). Rows
# that have no price yet keep the new cell blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column 22 = "V": inserting here pushes the current V ("nom") and W
# ("url_produit") one column to the right, to W and X respectively, and
# leaves a fresh, blank column V in their place.
$ws.Columns.Item(22).Insert()

# New header cell for the inserted timestamp column.
$ws.Range("V1").Value = "2026-01-28 15:21:05"

# Rows 2-80 hold a real price in column U; duplicate it into the new V
# column. Rows 81-205 have no price yet (U is blank), so V stays blank.
for ($row = 2; $row -le 80; $row++) {
    $price = $ws.Cells.Item($row, 21).Value()
    $ws.Cells.Item($row, 22).Value = $price
}
